$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.519.44"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "3.386.75"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.384.72"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "3.963.92"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "3.385.05"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.68%  "
$ws.Range("D18").Value = "61.603.89"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").Value = "3.527.67"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +8.76%  "
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("E28").Value = "  +4.88%  "
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("E31").Value = "  +5.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("B42").Value = "ONDO"
$ws.Range("C42").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.778"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("E44").Value = "  +8.44%  "
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.58%  "
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").Value = "2.346.85"
$ws.Range("E50").Value = "  +5.64%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0262"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.13%  "
